$wb = $excel.ActiveWorkbook

# --- Controls sheet: update saturating-DM input parameters ---
$controls = $wb.Worksheets.Item("Controls")
$controls.Range("B2").Value = 1000
$controls.Range("B3").Value = 31
$controls.Range("B4").Value = 100

# --- Make Controls the active/selected sheet (was Recruitment_Mortality) ---
$controls.Activate()
$controls.Range("B4").Select()

# --- Zoom the now-active window to 140% ---
$excel.ActiveWindow.Zoom = 140
